$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("computation_details")
$readme = $wb.Worksheets.Item("README")

# ---------------------------------------------------------------------------
# Register the "BATTERY" shared string before the longer README sentence
# that also mentions it, so the shared string table ends up in the same
# order as a natural edit (BATTERY rows added first, then the README blurb
# updated to mention it).
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "BATTERY"

# ---------------------------------------------------------------------------
# A1 header cell keeps its bold / wrap formatting but loses the explicit
# horizontal centering (only vertical centering + wrap remain). Build the
# new format on a scratch cell first (so it starts from a clean / default
# alignment) and paste the resulting format onto A1.
# ---------------------------------------------------------------------------
$scratch = $ws.Range("Z1")
$scratch.Font.Bold = $true
$scratch.WrapText = $true
$scratch.VerticalAlignment = -4108
$scratch.Copy()
$ws.Range("A1").PasteSpecial(-4122)
$scratch.Clear()

# ---------------------------------------------------------------------------
# Row 8: BATTERY / Total impact
# ---------------------------------------------------------------------------
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = "BATTERY"

$ws.Range("B5").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("B8").Value = "Total impact"

$ws.Range("C2").Copy()
$ws.Range("C8:H8").PasteSpecial(-4122)

$ws.Range("C8").Value = 11.1739394657824
$ws.Range("D8").Value = 20
$ws.Range("E8").Formula = "=C8/D8"
$ws.Range("F8").Value = 10
$ws.Range("G8").Formula = "=E8*F8"
$ws.Range("H8").Value = 30

$ws.Range("I8").Interior.Color = 65535
$ws.Range("I8").HorizontalAlignment = -4108
$ws.Range("I8").VerticalAlignment = -4108
$ws.Range("I8").Formula = "=G8*H8"

# ---------------------------------------------------------------------------
# Row 9: BATTERY / Contribution of Carbon dioxide, fossil
# ---------------------------------------------------------------------------
$ws.Range("A7").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A9").Value = "BATTERY"

$ws.Range("B7").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("B9").Value = "Contribution of Carbon dioxide, fossil ('air', 'non-urban air or from high stacks')"

$ws.Range("C2").Copy()
$ws.Range("C9:H9").PasteSpecial(-4122)

$ws.Range("C9").Value = 6.5855610967956304
$ws.Range("D9").Value = 20
$ws.Range("E9").Formula = "=C9/D9"
$ws.Range("F9").Value = 10
$ws.Range("G9").Formula = "=E9*F9"
$ws.Range("H9").Value = 30

$ws.Range("I8").Copy()
$ws.Range("I9").PasteSpecial(-4122)
$ws.Range("I9").Formula = "=G9*H9"

$ws.Rows.Item(9).RowHeight = 29

# ---------------------------------------------------------------------------
# Row 10: trailing empty formatted cell
# ---------------------------------------------------------------------------
$ws.Range("A7").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").ClearContents()

$ws.Activate()
$ws.Range("I10").Select()

# ---------------------------------------------------------------------------
# Sheet "README": update the intro paragraph to mention the new BATTERY
# technology, and move the active selection to A3.
# ---------------------------------------------------------------------------
$readme.Range("A2").Value = "We consider a single technology (BATTERY), and a composition of two sub-technologies, namely TRAIN_FREIGHT_DIESEL_LOC and TRAIN_FREIGHT_DIESEL_WAG, which are combined into the technology TRAIN_FREIGHT_DIESEL."

$readme.Activate()
$readme.Range("A3").Select()

# Re-activate computation_details so it remains the selected tab, matching
# the target workbook state (activeTab stays on computation_details).
$ws.Activate()
